$d = $word.ActiveDocument

# 1. Contract number
$d.Content.Find.Execute("CON6719273981243", $false, $false, $false, $false, $false, $true, 1, $false, "CON09262025", 2)

# 2. CC number + certification city (two separate changes in same run's text)
$d.Content.Find.Execute("con la CC No 1334567890, en su calidad", $false, $false, $false, $false, $false, $true, 1, $false, "con la CC No 1.029.384.756, en su calidad", 2)
$d.Content.Find.Execute("Cámara de Comercio de MEDELLÍN, la cual se adjunta", $false, $false, $false, $false, $false, $true, 1, $false, "Cámara de Comercio de MEDELLÍN, ANTIOQUIA, la cual se adjunta", 2)

# 3. "la ciudad MEDELLÍN." -> "la ciudad de MEDELLÍN."
$d.Content.Find.Execute("para contratos celebrados en la ciudad MEDELLÍN.", $false, $false, $false, $false, $false, $true, 1, $false, "para contratos celebrados en la ciudad de MEDELLÍN.", 2)

# 4. signing date
$d.Content.Find.Execute("el día ONCE (11) JULIO de 2025.", $false, $false, $false, $false, $false, $true, 1, $false, "el día VEINTISÉIS (26) SEPTIEMBRE de 2025.", 2)

# 5. C.C. No in signature block
$d.Content.Find.Execute("C.C. No 1334567890", $false, $false, $false, $false, $false, $true, 1, $false, "C.C. No 1.029.384.756", 2)

# 6. Company name in signature block
$d.Content.Find.Execute("IMB JUANCHO", $false, $false, $false, $false, $false, $true, 1, $false, "COMERCIALIZADORA EL POBLADO", 2)
